$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header section: account holder name and card/account number
$ws.Range("C2").Value = "Hartmut"

# B3 holds a 16-digit card number stored as TEXT (not a number). A direct
# .Value assignment of an all-digit string gets auto-coerced to a numeric
# type by the engine (like real Excel "smart" input), so stage the text in
# a scratch cell formatted as Text, then paste-special (values only) into
# B3. This carries over the string type without disturbing B3's own style.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "2570314725427075"
$ws.Range("Z1").Copy()
$ws.Range("B3").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("C3").Value = "Mohaupt"

# Opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 14.07.2025"

# Transaction row 6
$ws.Range("B6").Value = "16.07."
$ws.Range("C6").Value = "17.07."
$ws.Range("D6").Value = "EBAY MKTPLC EU BTACUT"
$ws.Range("E6").Value = "28,50-"

# Transaction row 7
$ws.Range("B7").Value = "18.07."
$ws.Range("C7").Value = "19.07."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-66425710"
$ws.Range("E7").Value = "53,43-"

# Transaction row 8
$ws.Range("B8").Value = "19.07."
$ws.Range("C8").Value = "20.07."
$ws.Range("D8").Value = "EBAY MKTPLC EU GOFPHQ"
$ws.Range("E8").Value = "53,32-"

# Transaction row 9 is removed (cleared), only 3 transactions remain now
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("E9").WrapText = $true

# Closing balance amount
$ws.Range("E12").Value = "135,25-"

# Next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 01.08.2025"
